$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row above row 5 ("History of partnering with Wannon Water"
#    / "Very important!") - everything from the old row 5 downward shifts
#    down by one row.
# ---------------------------------------------------------------------------
$ws.Rows("5:5").Insert()

# Copy the formatting (styles) from the row directly below (the old row 5,
# now row 6) onto the freshly inserted row 5, then restore the row height
# which PasteSpecial does not carry across.
$ws.Range("A6:C6").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)
$ws.Rows("5:5").RowHeight = 19.95

# Populate the new row's content.
$ws.Range("A5").Value = "History of partnering with Wannon Water"
$ws.Range("C5").Value = "Very important!"

# ---------------------------------------------------------------------------
# 2. The old row 5 (now row 6) comments cell gets shortened.
# ---------------------------------------------------------------------------
$ws.Range("C6").Value = "Higher score if company has previous worked with Wannon Water."

# ---------------------------------------------------------------------------
# 3. The "service will be delivered through the following channels" comment
#    (old row 11, now row 12) gets an extra closing sentence appended.
# ---------------------------------------------------------------------------
$deliveryText = @"
The service will be delivered through the following channels:
Research Reports: Access to detailed reports on trends, best practices, and technologies in IT and OT.
Market Analysis: Regular updates on the competitive landscape, including vendor evaluations, market forecasts, and technology adoption rates.
Strategic Guidance: Customised strategic advice based on the latest industry trends and specific business needs.
Benchmarking Services: Data and tools to compare Wannon Water's performance against industry peers.
Access to Analysts: Direct consultations with industry experts for personalised advice and support.
Workshops, Webinars, and Conferences: Access to world-class conferences, educational sessions, and workshops that provide cutting-edge insights and networking opportunities with industry leaders.
Specifically call out the number of seats or employee access licences as part of the proposal.
"@
$ws.Range("C12").Value = $deliveryText

# ---------------------------------------------------------------------------
# 4. Keep the sheet view / selection consistent with the edited area.
# ---------------------------------------------------------------------------
$ws.Range("C12").Select()
$excel.ActiveWindow.ScrollRow = 11
